# Apply updated crypto price/volume snapshot (GitHub Actions refresh, 2023-01-29 17:52 UTC).
# Source data reordered a handful of coins and refreshed their quotes; some other
# rows only received refreshed Price/Volume(1h) quotes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin / Link columns (plain text, never numeric-looking) ---
$nameLinkUpdates = @(
    @{ Cell = 'B8'; Value = 'MXToken' }
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' }
    @{ Cell = 'B9'; Value = 'LiechtensteinCryptoassetsExchange' }
    @{ Cell = 'C9'; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' }
    @{ Cell = 'B10'; Value = 'WazirX' }
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx' }
    @{ Cell = 'B11'; Value = 'MandalaExchangeToken' }
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' }
    @{ Cell = 'B12'; Value = 'BitrueCoin' }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' }
    @{ Cell = 'B13'; Value = 'BitMartToken' }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx' }
    @{ Cell = 'B14'; Value = 'BitForexToken' }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' }
    @{ Cell = 'B15'; Value = 'TigerCash' }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' }
    @{ Cell = 'B16'; Value = 'LEO' }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' }
    @{ Cell = 'B17'; Value = 'GateToken' }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt' }
    @{ Cell = 'B48'; Value = 'BOLO' }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo' }
    @{ Cell = 'B49'; Value = 'CoinbaseStockToken' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin' }
)
foreach ($u in $nameLinkUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# --- Price / Volume(1h) columns ---
# These look like numbers/percentages, so the cell is forced to Text format ("@")
# before assignment; otherwise Excel would silently convert them to numeric values
# and the exact original string formatting (trailing zeros, literal "%") would be lost.
$priceVolumeUpdates = @(
    @{ Cell = 'D2'; Value = '318.15' }
    @{ Cell = 'E2'; Value = '3.82%' }
    @{ Cell = 'D3'; Value = '39.80' }
    @{ Cell = 'E3'; Value = '2.60%' }
    @{ Cell = 'D4'; Value = '5.141' }
    @{ Cell = 'E4'; Value = '0.73%' }
    @{ Cell = 'D5'; Value = '0.08214' }
    @{ Cell = 'E5'; Value = '2.02%' }
    @{ Cell = 'D6'; Value = '2.134' }
    @{ Cell = 'E6'; Value = '10.02%' }
    @{ Cell = 'D7'; Value = '8.309' }
    @{ Cell = 'E7'; Value = '3.85%' }
    @{ Cell = 'D8'; Value = '0.9349' }
    @{ Cell = 'E8'; Value = '0.38%' }
    @{ Cell = 'D9'; Value = '0.1392' }
    @{ Cell = 'E9'; Value = '-4.37%' }
    @{ Cell = 'D10'; Value = '0.1986' }
    @{ Cell = 'E10'; Value = '3.22%' }
    @{ Cell = 'D11'; Value = '0.09090' }
    @{ Cell = 'E11'; Value = '1.08%' }
    @{ Cell = 'D12'; Value = '0.03476' }
    @{ Cell = 'E12'; Value = '-0.81%' }
    @{ Cell = 'D13'; Value = '0.09799' }
    @{ Cell = 'E13'; Value = '0.11%' }
    @{ Cell = 'D14'; Value = '0.001397' }
    @{ Cell = 'E14'; Value = '0.27%' }
    @{ Cell = 'D15'; Value = '0.006302' }
    @{ Cell = 'E15'; Value = '7.47%' }
    @{ Cell = 'D16'; Value = '3.675' }
    @{ Cell = 'E16'; Value = '-2.89%' }
    @{ Cell = 'D17'; Value = '4.285' }
    @{ Cell = 'E17'; Value = '2.07%' }
    @{ Cell = 'E18'; Value = '-2.68%' }
    @{ Cell = 'D19'; Value = '0.3474' }
    @{ Cell = 'E19'; Value = '1.58%' }
    @{ Cell = 'D20'; Value = '0.1291' }
    @{ Cell = 'E20'; Value = '-0.92%' }
    @{ Cell = 'D21'; Value = '4.898' }
    @{ Cell = 'E21'; Value = '2.47%' }
    @{ Cell = 'E22'; Value = '1.45%' }
    @{ Cell = 'D23'; Value = '0.04323' }
    @{ Cell = 'E23'; Value = '-1.47%' }
    @{ Cell = 'D24'; Value = '0.001225' }
    @{ Cell = 'E24'; Value = '-1.07%' }
    @{ Cell = 'D25'; Value = '0.004763' }
    @{ Cell = 'E25'; Value = '11.42%' }
    @{ Cell = 'D26'; Value = '0.0001299' }
    @{ Cell = 'E26'; Value = '-0.16%' }
    @{ Cell = 'D27'; Value = '0.0003996' }
    @{ Cell = 'E27'; Value = '-10.14%' }
    @{ Cell = 'D39'; Value = '0.02236' }
    @{ Cell = 'E39'; Value = '9.60%' }
    @{ Cell = 'D40'; Value = '0.05224' }
    @{ Cell = 'E40'; Value = '3.89%' }
    @{ Cell = 'D41'; Value = '0.007496' }
    @{ Cell = 'E41'; Value = '0.81%' }
    @{ Cell = 'D42'; Value = '0.009610' }
    @{ Cell = 'E42'; Value = '-3.80%' }
    @{ Cell = 'E43'; Value = '2.32%' }
    @{ Cell = 'D44'; Value = '0.002149' }
    @{ Cell = 'E44'; Value = '1.25%' }
    @{ Cell = 'D45'; Value = '0.009511' }
    @{ Cell = 'E45'; Value = '5.28%' }
    @{ Cell = 'D46'; Value = '0.00006607' }
    @{ Cell = 'E46'; Value = '6.58%' }
    @{ Cell = 'D47'; Value = '0.00000000749' }
    @{ Cell = 'E47'; Value = '-0.19%' }
    @{ Cell = 'D48'; Value = '0.002771' }
    @{ Cell = 'E48'; Value = '-0.47%' }
    @{ Cell = 'D49'; Value = '0.001200' }
    @{ Cell = 'E49'; Value = '-25.07%' }
    @{ Cell = 'D50'; Value = '0.00002098' }
    @{ Cell = 'E50'; Value = '-0.19%' }
    @{ Cell = 'D51'; Value = '0.0001998' }
    @{ Cell = 'E51'; Value = '-0.19%' }
)
foreach ($u in $priceVolumeUpdates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
